# "Clean ATDD sheet.xlsx" - replace the placeholder example content on the
# "ATDD Scenarios" sheet with a more meaningful Feature / Sub Feature name
# and meaningful Given/When/Then descriptions (instead of the bare letters
# A..I), then tidy up the resulting column widths / row heights the same
# way Excel does after the text changes.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ATDD Scenarios")
$ws2 = $wb.Worksheets.Item("Examples")

# --- Feature / Sub Feature (column A & B) for every data row ---------------
$ws1.Range("A2:A14").Value = "MyFeature"
$ws1.Range("B2:B14").Value = "MySubFeature"

# --- Given-When-Then (Description) column (G) -------------------------------
$ws1.Range("G4").Value  = "Some record"
$ws1.Range("G5").Value  = "Do something"
$ws1.Range("G6").Value  = "Something happens"
$ws1.Range("G8").Value  = "Another record"
$ws1.Range("G9").Value  = "Do something else"
$ws1.Range("G10").Value = "An error was thrown"
$ws1.Range("G12").Value = "Some other Record"
$ws1.Range("G13").Value = "Do something in a page"
$ws1.Range("G14").Value = "Something else happens"

# --- Column widths now that the text is longer/shorter than before ---------
$ws1.Columns.Item(5).ColumnWidth  = 13.0
$ws1.Columns.Item(7).ColumnWidth  = 22.5
$ws1.Columns.Item(9).ColumnWidth  = 29.0
$ws1.Columns.Item(10).ColumnWidth = 30.166666666666668
$ws1.Columns.Item(11).ColumnWidth = 33.333333333333336

# --- Row heights collapse back to the single-line default (15) -------------
$ws1.Rows.Item(2).RowHeight  = 15
$ws1.Rows.Item(4).RowHeight  = 15
$ws1.Rows.Item(5).RowHeight  = 15
$ws1.Rows.Item(6).RowHeight  = 15
$ws1.Rows.Item(7).RowHeight  = 15
$ws1.Rows.Item(8).RowHeight  = 15
$ws1.Rows.Item(9).RowHeight  = 15
$ws1.Rows.Item(10).RowHeight = 15
$ws1.Rows.Item(11).RowHeight = 15
$ws1.Rows.Item(12).RowHeight = 15
$ws1.Rows.Item(13).RowHeight = 15
$ws1.Rows.Item(14).RowHeight = 15

# --- Restore the on-screen selection: user ended up on "Examples"!L17 while
#     saving, then flipped back to "ATDD Scenarios" which stays the active
#     (visible) tab with its cursor back at A1. -----------------------------
$ws2.Activate()
$ws2.Range("L17").Select()
$ws1.Activate()
$ws1.Range("A1").Select()
